$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "[-, -, -, -]"

$ws.Range("C3").Value = "[-, 'MEC-3B-Metrologia 2', -, -]"
$ws.Range("D3").Value = "-"

$ws.Range("C4").Value = "[-, 'MEC-3B-Metrologia 2', -, -]"
$ws.Range("D4").Value = "-"

$ws.Range("C6").Value = "[-, 'MEC-3B-Metrologia 2', -, -]"
$ws.Range("D6").Value = "-"

$ws.Range("C7").Value = "[Ismail-Metrologia 1-1A, 'MEC-3B-Metrologia 2', -, -]"
$ws.Range("D7").Value = "-"

$ws.Range("B10").Value = "-"
$ws.Range("D10").Value = "[-, -, -, -]"

$ws.Range("B11").Value = "[Ismail-Metrologia 1-1A, -, -, -]"

$ws.Range("B12").Value = "[Ismail-Metrologia 1-1A, -, -, -]"

$ws.Range("B14").Value = "[Ismail-Metrologia 1-1A, -, -, -]"

$ws.Range("B15").Value = "['MEC-1A-Metrologia 1', -, -, -]"

$ws.Range("B16").Value = "-"

$ws.Range("B18").Value = "[-, -, 'MEC-2NA-Metrologia 2', -]"
$ws.Range("D18").Value = "[-, -, -, -]"
$ws.Range("E18").Value = "-"

$ws.Range("B19").Value = "[-, -, 'MEC-2NA-Metrologia 2', -]"
$ws.Range("E19").Value = "-"

$ws.Range("B20").Value = "[-, -, 'MEC-2NA-Metrologia 2', -]"
$ws.Range("C20").Value = "[-, -, 'MEC-2NB-Metrologia 2', -]"
$ws.Range("D20").Value = "[-, -, -, -]"

$ws.Range("B21").Value = "[-, -, 'MEC-2NA-Metrologia 2', -]"
$ws.Range("C21").Value = "[-, -, Leandro-M. Maq. E. I.-2NB, -]"
$ws.Range("E21").Value = "[-, -, -, -]"
